# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets,
# matching the regenerated gh-pages data snapshot (commit 456a3b4).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 47
$ws1.Range("F6").Value = 27
$ws1.Range("F7").Value = 576
$ws1.Range("F8").Value = 51
$ws1.Range("F9").Value = 8362
$ws1.Range("F10").Value = 778
$ws1.Range("F11").Value = 306
$ws1.Range("F12").Value = 1120
$ws1.Range("F13").Value = 879
$ws1.Range("F14").Value = 65
$ws1.Range("F16").Value = 218
$ws1.Range("F17").Value = 145
$ws1.Range("F19").Value = 217
$ws1.Range("F20").Value = 922

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 47
$ws4.Range("F7").Value = 27
$ws4.Range("F9").Value = 576
$ws4.Range("F10").Value = 51
$ws4.Range("F11").Value = 8362
$ws4.Range("F12").Value = 778
$ws4.Range("F13").Value = 306
$ws4.Range("F14").Value = 1120
$ws4.Range("F15").Value = 879
$ws4.Range("F16").Value = 65
$ws4.Range("F18").Value = 218
$ws4.Range("F19").Value = 145
$ws4.Range("F21").Value = 217
$ws4.Range("F22").Value = 922

$wb.Save()
